$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1658291457286432
$ws.Range("C2").Value = 0.628140703517588
$ws.Range("J2").Value = 0.01005025125628141
$ws.Range("P2").Value = 0.1407035175879397
$ws.Range("S2").Value = 0.05527638190954774

# Row 3
$ws.Range("B3").Value = 0.01587301587301587
$ws.Range("J3").Value = 0.02380952380952381
$ws.Range("P3").Value = 0.8015873015873016
$ws.Range("S3").Value = 0.1587301587301587

# Row 4
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3

# Row 6
$ws.Range("B6").Value = 0.04265402843601896
$ws.Range("D6").Value = 0.009478672985781991
$ws.Range("F6").Value = 0.07109004739336493
$ws.Range("J6").Value = 0.2464454976303317
$ws.Range("O6").Value = 0.009478672985781991
$ws.Range("Q6").Value = 0.1706161137440758
$ws.Range("R6").Value = 0.06635071090047394
$ws.Range("S6").Value = 0.3838862559241706

# Row 7
$ws.Range("B7").Value = 0.05617977528089887
$ws.Range("D7").Value = 0.02247191011235955
$ws.Range("E7").Value = 0.01123595505617977
$ws.Range("F7").Value = 0.06179775280898876
$ws.Range("J7").Value = 0.1460674157303371
$ws.Range("O7").Value = 0.01685393258426966
$ws.Range("Q7").Value = 0.2078651685393259
$ws.Range("R7").Value = 0.06179775280898876
$ws.Range("S7").Value = 0.4157303370786517

# Row 8
$ws.Range("B8").Value = 0.06108597285067873
$ws.Range("D8").Value = 0.02036199095022624
$ws.Range("F8").Value = 0.05429864253393665
$ws.Range("J8").Value = 0.1063348416289593
$ws.Range("O8").Value = 0.02262443438914027
$ws.Range("Q8").Value = 0.1877828054298643
$ws.Range("R8").Value = 0.08597285067873303
$ws.Range("S8").Value = 0.4615384615384616

# Row 9
$ws.Range("B9").Value = 0.08823529411764706
$ws.Range("D9").Value = 0.004201680672268907
$ws.Range("F9").Value = 0.06302521008403361
$ws.Range("J9").Value = 0.07142857142857142
$ws.Range("O9").Value = 0.01680672268907563
$ws.Range("Q9").Value = 0.2521008403361344
$ws.Range("R9").Value = 0.09243697478991597
$ws.Range("S9").Value = 0.4117647058823529

# Row 10
$ws.Range("B10").Value = 0.08541846419327007
$ws.Range("D10").Value = 0.01294219154443486
$ws.Range("E10").Value = 0.0008628127696289905
$ws.Range("F10").Value = 0.07506471095772217
$ws.Range("J10").Value = 0.091458153580673
$ws.Range("O10").Value = 0.01553062985332183
$ws.Range("Q10").Value = 0.2018981880931838
$ws.Range("R10").Value = 0.1087144089732528
$ws.Range("S10").Value = 0.4081104400345125

# Row 11
$ws.Range("G11").Value = 0.1721611721611722
$ws.Range("J11").Value = 0.06593406593406594
$ws.Range("K11").Value = 0.2234432234432235
$ws.Range("L11").Value = 0.5347985347985348
$ws.Range("S11").Value = 0.003663003663003663

# Row 12
$ws.Range("G12").Value = 0.6794871794871795
$ws.Range("J12").Value = 0.1987179487179487
$ws.Range("K12").Value = 0.03205128205128205
$ws.Range("L12").Value = 0.04487179487179487
$ws.Range("S12").Value = 0.04487179487179487

# Row 13
$ws.Range("G13").Value = 0.7857142857142857
$ws.Range("J13").Value = 0.1666666666666667
$ws.Range("S13").Value = 0.04761904761904762

# Row 15
$ws.Range("F15").Value = 0.01395348837209302
$ws.Range("H15").Value = 0.1395348837209302
$ws.Range("I15").Value = 0.09302325581395349
$ws.Range("J15").Value = 0.3581395348837209
$ws.Range("K15").Value = 0.04186046511627907
$ws.Range("M15").Value = 0.0186046511627907
$ws.Range("N15").Value = 0.009302325581395349
$ws.Range("O15").Value = 0.04651162790697674
$ws.Range("S15").Value = 0.2790697674418605

# Row 16
$ws.Range("F16").Value = 0.01351351351351351
$ws.Range("H16").Value = 0.1824324324324324
$ws.Range("I16").Value = 0.1013513513513514
$ws.Range("J16").Value = 0.3851351351351351
$ws.Range("K16").Value = 0.1216216216216216
$ws.Range("M16").Value = 0.01351351351351351
$ws.Range("O16").Value = 0.02027027027027027
$ws.Range("S16").Value = 0.1621621621621622

# Row 17
$ws.Range("F17").Value = 0.0244988864142539
$ws.Range("H17").Value = 0.178173719376392
$ws.Range("I17").Value = 0.1135857461024499
$ws.Range("J17").Value = 0.3942093541202673
$ws.Range("K17").Value = 0.08908685968819599
$ws.Range("M17").Value = 0.0178173719376392
$ws.Range("N17").Value = 0.0022271714922049
$ws.Range("O17").Value = 0.06904231625835189
$ws.Range("S17").Value = 0.111358574610245

# Row 18
$ws.Range("F18").Value = 0.01442307692307692
$ws.Range("H18").Value = 0.1778846153846154
$ws.Range("I18").Value = 0.1298076923076923
$ws.Range("J18").Value = 0.4230769230769231
$ws.Range("K18").Value = 0.0625
$ws.Range("M18").Value = 0.01442307692307692
$ws.Range("N18").Value = 0.004807692307692308
$ws.Range("O18").Value = 0.07211538461538461
$ws.Range("S18").Value = 0.1009615384615385

# Row 19
$ws.Range("F19").Value = 0.01422924901185771
$ws.Range("H19").Value = 0.2150197628458498
$ws.Range("I19").Value = 0.1003952569169961
$ws.Range("J19").Value = 0.3683794466403162
$ws.Range("K19").Value = 0.09723320158102766
$ws.Range("M19").Value = 0.01976284584980237
$ws.Range("N19").Value = 0.0007905138339920949
$ws.Range("O19").Value = 0.07114624505928854
$ws.Range("S19").Value = 0.1130434782608696
